$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell (H1) onto the new headers
# so they get the same style (bold, border, centered) as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# I and J values for rows 2-31
$values = @(
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(4, 5),
    @(5, 5),
    @(5, 5),
    @(7, 7),
    @(5, 5),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
